# Update product listing rows 30-75 of "Resultados da pesquisa":
#  - rows 30-37 become new Pingo Doce egg/ovo entries (replacing the old
#    "Arroz Agulha" rows that used to sit there)
#  - the rice ("Arroz ...") rows that used to be at 30-42 shift down to 38-51
#    (two of them also get an updated price)
#  - rows 52-75 are brand-new Continente egg / clara de ovo entries
# Net effect: used range grows from A1:C42 to A1:C75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(30, 1).Value = "Ovos de Galinhas Criadas ao Ar Livre Classe M/L Pingo Doce 6 un"
$ws.Cells.Item(30, 2).Value = "1,59€ / un"
$ws.Cells.Item(30, 3).Value = "Pingo Doce"

$ws.Cells.Item(31, 1).Value = "Ovos de Codorniz Pingo Doce 12 un"
$ws.Cells.Item(31, 2).Value = "1,30€ / un"
$ws.Cells.Item(31, 3).Value = "Pingo Doce"

$ws.Cells.Item(32, 1).Value = "Ovos de Galinhas Criadas no Solo Classe L Pingo Doce 6 un"
$ws.Cells.Item(32, 2).Value = "1,48€ / un"
$ws.Cells.Item(32, 3).Value = "Pingo Doce"

$ws.Cells.Item(33, 1).Value = "Ovos de Galinhas Criadas no Solo Classe L Pingo Doce 12 un"
$ws.Cells.Item(33, 2).Value = "2,48€ / un"
$ws.Cells.Item(33, 3).Value = "Pingo Doce"

$ws.Cells.Item(34, 1).Value = "Ovos de Galinhas Criadas no Solo Classe M Pingo Doce 6 un"
$ws.Cells.Item(34, 2).Value = "1,38€ / un"
$ws.Cells.Item(34, 3).Value = "Pingo Doce"

$ws.Cells.Item(35, 1).Value = "Ovos de Galinhas Criadas no Solo Classe M Pingo Doce 12 un"
$ws.Cells.Item(35, 2).Value = "2,33€ / un"
$ws.Cells.Item(35, 3).Value = "Pingo Doce"

$ws.Cells.Item(36, 1).Value = "Ovos de Galinhas Criadas no Solo Classe XL Pingo Doce 6 un"
$ws.Cells.Item(36, 2).Value = "2,07€ / un"
$ws.Cells.Item(36, 3).Value = "Pingo Doce"

$ws.Cells.Item(37, 1).Value = "Ovos Classe M/L Pingo Doce Biológico 6 un"
$ws.Cells.Item(37, 2).Value = "2,39€ / un"
$ws.Cells.Item(37, 3).Value = "Pingo Doce"

$ws.Cells.Item(38, 1).Value = "Arroz Agulha"
$ws.Cells.Item(38, 2).Value = "€1,69`n/un"
$ws.Cells.Item(38, 3).Value = "Continente"

$ws.Cells.Item(39, 1).Value = "Arroz Agulha"
$ws.Cells.Item(39, 2).Value = "€1,59`n/un"
$ws.Cells.Item(39, 3).Value = "Continente"

$ws.Cells.Item(40, 1).Value = "Arroz Agulha"
$ws.Cells.Item(40, 2).Value = "€1,33`n/un"
$ws.Cells.Item(40, 3).Value = "Continente"

$ws.Cells.Item(41, 1).Value = "Arroz Agulha"
$ws.Cells.Item(41, 2).Value = "€1,38`n/un"
$ws.Cells.Item(41, 3).Value = "Continente"

$ws.Cells.Item(42, 1).Value = "Arroz Agulha Europa"
$ws.Cells.Item(42, 2).Value = "€1,18`n/un"
$ws.Cells.Item(42, 3).Value = "Continente"

$ws.Cells.Item(43, 1).Value = "Arroz Agulha Selecionado"
$ws.Cells.Item(43, 2).Value = "€1,26`n/un"
$ws.Cells.Item(43, 3).Value = "Continente"

$ws.Cells.Item(44, 1).Value = "Arroz Agulha Integral"
$ws.Cells.Item(44, 2).Value = "€1,28`n/un"
$ws.Cells.Item(44, 3).Value = "Continente"

$ws.Cells.Item(45, 1).Value = "Arroz Agulha Extra Longo"
$ws.Cells.Item(45, 2).Value = "€2,37`n/un"
$ws.Cells.Item(45, 3).Value = "Continente"

$ws.Cells.Item(46, 1).Value = "Arroz Agulha Extra Longo"
$ws.Cells.Item(46, 2).Value = "€1,28`n/un"
$ws.Cells.Item(46, 3).Value = "Continente"

$ws.Cells.Item(47, 1).Value = "Arroz Agulha Europa Pack Poupança"
$ws.Cells.Item(47, 2).Value = "€1,36`n/un"
$ws.Cells.Item(47, 3).Value = "Continente"

$ws.Cells.Item(48, 1).Value = "Arroz Agulha Branqueado Extra Longo"
$ws.Cells.Item(48, 2).Value = "€1,33`n/un"
$ws.Cells.Item(48, 3).Value = "Continente"

$ws.Cells.Item(49, 1).Value = "Arroz Agulha Pronto a Comer"
$ws.Cells.Item(49, 2).Value = "€1,45`n/un"
$ws.Cells.Item(49, 3).Value = "Continente"

$ws.Cells.Item(50, 1).Value = "Arroz Agulha Pronto a Comer sem Glúten"
$ws.Cells.Item(50, 2).Value = "€1,45`n/un"
$ws.Cells.Item(50, 3).Value = "Continente"

$ws.Cells.Item(51, 1).Value = "Arroz Longo Comum"
$ws.Cells.Item(51, 2).Value = "€1,08`n/un"
$ws.Cells.Item(51, 3).Value = "Continente"

$ws.Cells.Item(52, 1).Value = "Ovos de Solo Classe L"
$ws.Cells.Item(52, 2).Value = "€2,48`n/un"
$ws.Cells.Item(52, 3).Value = "Continente"

$ws.Cells.Item(53, 1).Value = "Ovos de Solo Classe M"
$ws.Cells.Item(53, 2).Value = "€2,33`n/un"
$ws.Cells.Item(53, 3).Value = "Continente"

$ws.Cells.Item(54, 1).Value = "Ovos de Ar Livre Classe M/L"
$ws.Cells.Item(54, 2).Value = "€3,74`n/un"
$ws.Cells.Item(54, 3).Value = "Continente"

$ws.Cells.Item(55, 1).Value = "Ovos Classe M"
$ws.Cells.Item(55, 2).Value = "€4,68`n/un"
$ws.Cells.Item(55, 3).Value = "Continente"

$ws.Cells.Item(56, 1).Value = "Ovos de Ar Livre Classe M/L"
$ws.Cells.Item(56, 2).Value = "€2,59`n/un"
$ws.Cells.Item(56, 3).Value = "Continente"

$ws.Cells.Item(57, 1).Value = "Ovos de Ar Livre"
$ws.Cells.Item(57, 2).Value = "€2,38`n/un"
$ws.Cells.Item(57, 3).Value = "Continente"

$ws.Cells.Item(58, 1).Value = "Ovos de Solo Classe XL"
$ws.Cells.Item(58, 2).Value = "€2,06`n/un"
$ws.Cells.Item(58, 3).Value = "Continente"

$ws.Cells.Item(59, 1).Value = "Ovos de Codorniz"
$ws.Cells.Item(59, 2).Value = "€1,30`n/un"
$ws.Cells.Item(59, 3).Value = "Continente"

$ws.Cells.Item(60, 1).Value = "Ovos de Solo Classe M"
$ws.Cells.Item(60, 2).Value = "€1,37`n/un"
$ws.Cells.Item(60, 3).Value = "Continente"

$ws.Cells.Item(61, 1).Value = "Ovos de Solo Classe L"
$ws.Cells.Item(61, 2).Value = "€1,48`n/un"
$ws.Cells.Item(61, 3).Value = "Continente"

$ws.Cells.Item(62, 1).Value = "Ovos Classe M/L"
$ws.Cells.Item(62, 2).Value = "€3,04`n/un"
$ws.Cells.Item(62, 3).Value = "Continente"

$ws.Cells.Item(63, 1).Value = "Clara de Ovo Pasteurizada"
$ws.Cells.Item(63, 2).Value = "€2,95`n/un"
$ws.Cells.Item(63, 3).Value = "Continente"

$ws.Cells.Item(64, 1).Value = "Ovos de Ar Livre Classe M/L"
$ws.Cells.Item(64, 2).Value = "€1,84`n/un"
$ws.Cells.Item(64, 3).Value = "Continente"

$ws.Cells.Item(65, 1).Value = "Ovos de Ar Livre Classe M/L"
$ws.Cells.Item(65, 2).Value = "€1,59`n/un"
$ws.Cells.Item(65, 3).Value = "Continente"

$ws.Cells.Item(66, 1).Value = "Ovos de Solo Classe M/L"
$ws.Cells.Item(66, 2).Value = "€3,94`n/un"
$ws.Cells.Item(66, 3).Value = "Continente"

$ws.Cells.Item(67, 1).Value = "Ovos Classe L"
$ws.Cells.Item(67, 2).Value = "€3,14`n/un"
$ws.Cells.Item(67, 3).Value = "Continente"

$ws.Cells.Item(68, 1).Value = "Clara de Ovo Líquida Pasteurizada"
$ws.Cells.Item(68, 2).Value = "€1,15`n/un"
$ws.Cells.Item(68, 3).Value = "Continente"

$ws.Cells.Item(69, 1).Value = "Ovos Classe M"
$ws.Cells.Item(69, 2).Value = "€3,04`n/un"
$ws.Cells.Item(69, 3).Value = "Continente"

$ws.Cells.Item(70, 1).Value = "Ovos Classe M/L Ruby"
$ws.Cells.Item(70, 2).Value = "€2,65`n/un"
$ws.Cells.Item(70, 3).Value = "Continente"

$ws.Cells.Item(71, 1).Value = "Ovos de Solo Classe S"
$ws.Cells.Item(71, 2).Value = "€1,18`n/un"
$ws.Cells.Item(71, 3).Value = "Continente"

$ws.Cells.Item(72, 1).Value = "Ovos de Ar Livre Classe XL"
$ws.Cells.Item(72, 2).Value = "€2,58`n/un"
$ws.Cells.Item(72, 3).Value = "Continente"

$ws.Cells.Item(73, 1).Value = "Ovos de Solo Classe M"
$ws.Cells.Item(73, 2).Value = "€3,15`n/un"
$ws.Cells.Item(73, 3).Value = "Continente"

$ws.Cells.Item(74, 1).Value = "Clara de Ovo Proteína Baunilha"
$ws.Cells.Item(74, 2).Value = "€1,84`n/un"
$ws.Cells.Item(74, 3).Value = "Continente"

$ws.Cells.Item(75, 1).Value = "Clara de Ovo Proteína Morango"
$ws.Cells.Item(75, 2).Value = "€1,84`n/un"
$ws.Cells.Item(75, 3).Value = "Continente"
